# Auto-generated edit script: updates computed market-price / profit columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 659.5
$ws.Range("I49").Value = 300
$ws.Range("J49").Value = 1019
$ws.Range("K49").Value = 900
$ws.Range("L49").Value = 3057
$ws.Range("M49").Value = -764
$ws.Range("N49").Value = -3329
$ws.Range("H106").Value = 2194.8333
$ws.Range("I106").Value = 2194.8333
$ws.Range("K106").Value = 2194.8333
$ws.Range("M106").Value = -1563.8333
$ws.Range("H121").Value = 1950
$ws.Range("J121").Value = 1900
$ws.Range("L121").Value = 5700
$ws.Range("N121").Value = -9194
$ws.Range("H132").Value = 6103246.5
$ws.Range("I132").Value = 8071358
$ws.Range("J132").Value = 2100.1
$ws.Range("K132").Value = 24214074
$ws.Range("L132").Value = 6300.299999999999
$ws.Range("M132").Value = -24211544
$ws.Range("N132").Value = -11360.3

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1247.58
$ws.Range("I32").Value = 966.57733
$ws.Range("J32").Value = 10333.333
$ws.Range("K32").Value = 966.57733
$ws.Range("L32").Value = 10333.333
$ws.Range("M32").Value = -679.57733
$ws.Range("N32").Value = -10907.333
$ws.Range("H61").Value = 1138.931
$ws.Range("I61").Value = 885.6957
$ws.Range("J61").Value = 2109.6667
$ws.Range("K61").Value = 885.6957
$ws.Range("L61").Value = 2109.6667
$ws.Range("M61").Value = -673.6957
$ws.Range("N61").Value = -2533.6667
$ws.Range("H74").Value = 4765.8
$ws.Range("J74").Value = 6482.143
$ws.Range("L74").Value = 6482.143
$ws.Range("N74").Value = -8230.143
$ws.Range("H77").Value = 4765.8
$ws.Range("J77").Value = 6482.143
$ws.Range("L77").Value = 32410.715
$ws.Range("N77").Value = -41146.715
$ws.Range("H88").Value = 1806.0588
$ws.Range("J88").Value = 1897.8334
$ws.Range("L88").Value = 1897.8334
$ws.Range("N88").Value = -2709.8334
$ws.Range("H91").Value = 1806.0588
$ws.Range("J91").Value = 1897.8334
$ws.Range("L91").Value = 1897.8334
$ws.Range("N91").Value = -4705.8334
$ws.Range("H132").Value = 3909.6943
$ws.Range("I132").Value = 4347.1113
$ws.Range("J132").Value = 2597.4443
$ws.Range("K132").Value = 13041.3339
$ws.Range("L132").Value = 7792.3329
$ws.Range("M132").Value = -10511.3339
$ws.Range("N132").Value = -12852.3329
$ws.Range("H136").Value = 1138.931
$ws.Range("I136").Value = 885.6957
$ws.Range("J136").Value = 2109.6667
$ws.Range("K136").Value = 2657.0871
$ws.Range("L136").Value = 6329.000100000001
$ws.Range("M136").Value = -107.0870999999997
$ws.Range("N136").Value = -11429.0001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 15986.8
$ws.Range("J35").Value = 15986.8
$ws.Range("L35").Value = 15986.8
$ws.Range("N35").Value = -16606.8
$ws.Range("H86").Value = 47323.918
$ws.Range("I86").Value = 70038.44
$ws.Range("K86").Value = 70038.44
$ws.Range("M86").Value = -68915.44
$ws.Range("H89").Value = 47323.918
$ws.Range("I89").Value = 70038.44
$ws.Range("K89").Value = 350192.2
$ws.Range("M89").Value = -344576.2
$ws.Range("H134").Value = 2732.56
$ws.Range("I134").Value = 2792.0908
$ws.Range("J134").Value = 2296
$ws.Range("K134").Value = 8376.2724
$ws.Range("L134").Value = 6888
$ws.Range("M134").Value = -5841.2724
$ws.Range("N134").Value = -11958

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H6").Value = 1089.091
$ws.Range("I6").Value = 1191.4286
$ws.Range("J6").Value = 910
$ws.Range("K6").Value = 1191.4286
$ws.Range("L6").Value = 910
$ws.Range("M6").Value = -1078.4286
$ws.Range("N6").Value = -1136
$ws.Range("H22").Value = 954.9167
$ws.Range("I22").Value = 443
$ws.Range("J22").Value = 1320.5714
$ws.Range("K22").Value = 443
$ws.Range("L22").Value = 1320.5714
$ws.Range("M22").Value = -93
$ws.Range("N22").Value = -2020.5714
$ws.Range("H31").Value = 20720.246
$ws.Range("I31").Value = 26184.35
$ws.Range("J31").Value = 3907.6155
$ws.Range("K31").Value = 26184.35
$ws.Range("L31").Value = 3907.6155
$ws.Range("M31").Value = -25889.35
$ws.Range("N31").Value = -4497.6155
$ws.Range("H34").Value = 20720.246
$ws.Range("I34").Value = 26184.35
$ws.Range("J34").Value = 3907.6155
$ws.Range("K34").Value = 26184.35
$ws.Range("L34").Value = 3907.6155
$ws.Range("M34").Value = -25982.35
$ws.Range("N34").Value = -4311.6155
$ws.Range("H45").Value = 12999.75
$ws.Range("I45").Value = 12333
$ws.Range("K45").Value = 12333
$ws.Range("M45").Value = -11740
$ws.Range("H99").Value = 8892.75
$ws.Range("I99").Value = 2909
$ws.Range("K99").Value = 2909
$ws.Range("M99").Value = -1411
$ws.Range("H105").Value = 1153.5883
$ws.Range("I105").Value = 1108.8889
$ws.Range("J105").Value = 1203.875
$ws.Range("K105").Value = 1108.8889
$ws.Range("L105").Value = 1203.875
$ws.Range("M105").Value = 638.1111000000001
$ws.Range("N105").Value = -4697.875
$ws.Range("H107").Value = 568.6
$ws.Range("I107").Value = 625.06665
$ws.Range("J107").Value = 483.9
$ws.Range("K107").Value = 625.06665
$ws.Range("L107").Value = 483.9
$ws.Range("M107").Value = 1294.93335
$ws.Range("N107").Value = -4323.9
$ws.Range("H122").Value = 11416.25
$ws.Range("J122").Value = 9666.333000000001
$ws.Range("L122").Value = 28998.999
$ws.Range("N122").Value = -33898.999
$ws.Range("H126").Value = 8892.75
$ws.Range("I126").Value = 2909
$ws.Range("K126").Value = 8727
$ws.Range("M126").Value = -6257
$ws.Range("H134").Value = 1910.9333
$ws.Range("I134").Value = 2261.111
$ws.Range("J134").Value = 1385.6666
$ws.Range("K134").Value = 6783.333
$ws.Range("L134").Value = 4156.9998
$ws.Range("M134").Value = -4248.333
$ws.Range("N134").Value = -9226.9998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1000000000
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H41").Value = 1850
$ws.Range("J41").Value = 3200
$ws.Range("L41").Value = 9600
$ws.Range("N41").Value = -10276
$ws.Range("H64").Value = 1000
$ws.Range("J64").Value = 1000
$ws.Range("L64").Value = 3000
$ws.Range("N64").Value = -3540
$ws.Range("H67").Value = 1000
$ws.Range("J67").Value = 1000
$ws.Range("L67").Value = 3000
$ws.Range("N67").Value = -4872
$ws.Range("H122").Value = 4163.222
$ws.Range("I122").Value = 304.73914
$ws.Range("K122").Value = 2742.65226
$ws.Range("M122").Value = -292.6522600000003
$ws.Range("H131").Value = 777.4
$ws.Range("J131").Value = 863.8674999999999
$ws.Range("L131").Value = 2591.6025
$ws.Range("N131").Value = -12671.6025

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 213.33333
$ws.Range("I13").Value = 213.33333
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 213.33333
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -74.33332999999999
$ws.Range("N13").ClearContents()
$ws.Range("H70").Value = 57395.74
$ws.Range("I70").Value = 91340.64999999999
$ws.Range("J70").Value = 5346.8667
$ws.Range("K70").Value = 91340.64999999999
$ws.Range("L70").Value = 5346.8667
$ws.Range("M70").Value = -91070.64999999999
$ws.Range("N70").Value = -5886.8667
$ws.Range("H73").Value = 57395.74
$ws.Range("I73").Value = 91340.64999999999
$ws.Range("J73").Value = 5346.8667
$ws.Range("K73").Value = 91340.64999999999
$ws.Range("L73").Value = 5346.8667
$ws.Range("M73").Value = -90404.64999999999
$ws.Range("N73").Value = -7218.8667
$ws.Range("H132").Value = 2391.5
$ws.Range("I132").Value = 1791
$ws.Range("J132").Value = 3352.3
$ws.Range("K132").Value = 5373
$ws.Range("L132").Value = 10056.9
$ws.Range("M132").Value = -2843
$ws.Range("N132").Value = -15116.9

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4202.2
$ws.Range("I132").Value = 4252.5625
$ws.Range("J132").Value = 4000.75
$ws.Range("K132").Value = 12757.6875
$ws.Range("L132").Value = 12002.25
$ws.Range("M132").Value = -10227.6875
$ws.Range("N132").Value = -17062.25
$ws.Range("H136").Value = 1815.3846
$ws.Range("I136").Value = 1783.3334
$ws.Range("J136").Value = 1842.8572
$ws.Range("K136").Value = 5350.0002
$ws.Range("L136").Value = 5528.571599999999
$ws.Range("M136").Value = -2800.0002
$ws.Range("N136").Value = -10628.5716

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1951.2222
$ws.Range("I126").Value = 2008.762
$ws.Range("J126").Value = 1749.8334
$ws.Range("K126").Value = 6026.286
$ws.Range("L126").Value = 5249.5002
$ws.Range("M126").Value = -3556.286
$ws.Range("N126").Value = -10189.5002
$ws.Range("H132").Value = 3362.0454
$ws.Range("I132").Value = 3268.647
$ws.Range("J132").Value = 3679.6
$ws.Range("K132").Value = 9805.940999999999
$ws.Range("L132").Value = 11038.8
$ws.Range("M132").Value = -7275.940999999999
$ws.Range("N132").Value = -16098.8
